$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.645.57'
$ws.Range("E2").Value = '  +0.95%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.564.48'
$ws.Range("E3").Value = '  -0.06%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.21%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.38'
$ws.Range("E5").Value = '  -0.30%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.513'
$ws.Range("E6").Value = '  +4.82%  '

# Row 7
$ws.Range("E7").Value = '  -0.25%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.86'
$ws.Range("E8").Value = '  +5.33%  '

# Row 9
$ws.Range("E9").Value = '  +0.80%  '

# Row 10
$ws.Range("E10").Value = '  -0.23%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0899'
$ws.Range("E11").Value = '  +0.64%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.787.16'

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.577.11'
$ws.Range("E13").Value = '  +0.75%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.676.81'
$ws.Range("E14").Value = '  +1.12%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.517'
$ws.Range("E15").Value = '  +0.69%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.64'
$ws.Range("E16").Value = '  -0.84%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.60'
$ws.Range("E17").Value = '  +1.76%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '227.79'
$ws.Range("E18").Value = '  -0.19%  '

# Row 19
$ws.Range("E19").Value = '  -0.96%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0682'
$ws.Range("E20").Value = '  +0.31%  '

# Row 21
$ws.Range("E21").Value = '  -0.18%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.93'
$ws.Range("E22").Value = '  -0.21%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.08'
$ws.Range("E23").Value = '  +1.49%  '

# Row 24
$ws.Range("E24").Value = '  +0.67%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.07'
$ws.Range("E25").Value = '  +1.16%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.107'
$ws.Range("E26").Value = '  +3.24%  '

# Row 27
$ws.Range("E27").Value = '  -0.70%  '

# Row 28
$ws.Range("E28").Value = '  -0.21%  '

# Row 29
$ws.Range("E29").Value = '  -1.07%  '

# Row 30
$ws.Range("E30").Value = '  -3.81%  '

# Row 31
$ws.Range("E31").Value = '  -0.64%  '

# Row 32
$ws.Range("E32").Value = '  +0.07%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.402.30'
$ws.Range("E33").Value = '  +1.20%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.01'
$ws.Range("E34").Value = '  -2.42%  '

# Row 35
$ws.Range("E35").Value = '  -3.45%  '

# Row 36
$ws.Range("E36").Value = '  -2.11%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.69'
$ws.Range("E37").Value = '  +1.54%  '

# Row 38
$ws.Range("E38").Value = '  -2.08%  '

# Row 39
$ws.Range("E39").Value = '  +0.19%  '

# Row 40
$ws.Range("B40").Value = 'ImmutableX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.518'

# Row 41
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.94'
$ws.Range("E41").Value = '  -0.28%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  -0.18%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0463'
$ws.Range("E43").Value = '  -1.25%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.767'
$ws.Range("E44").Value = '  -2.27%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.98'
$ws.Range("E45").Value = '  +2.79%  '

# Row 46
$ws.Range("E46").Value = '  -2.42%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.699.82'
$ws.Range("E47").Value = '  -0.05%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.849'
$ws.Range("E48").Value = '  -7.75%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '84.76'
$ws.Range("E49").Value = '  -0.85%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '42.61'
$ws.Range("E50").Value = '  +4.82%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0512'
$ws.Range("E51").Value = '  -0.48%  '
